# Adapt DE results table: rename "Gen" header to "MaxFES", switch the
# x-axis values from raw generation/FES counts to normalised fractions,
# drop the "Run 50" run (column) and recompute the trailing "Mean" column
# over the remaining 50 runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header rename: A1 "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2) New normalised values for column A (rows 2-14)
$genValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $genValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $genValues[$i]
}

# 3) Remove the "Run 50" column entirely (was column AZ); this shifts the
#    trailing "Mean" column left from BA to AZ automatically.
$ws.Range("AZ1:AZ14").EntireColumn.Delete()

# 4) Recompute the "Mean" column (now AZ) over the remaining 50 runs
#    (columns B:AY) for each data row, rounded to 8 decimal places to match
#    the precision used throughout the rest of the table.
$means = @(237.11098024, 224.16799898, 115.61130607, 6.24205313, 1.09589465, 0.70192377, 0.6443137, 0.61342089, 0.55786512, 0.49684318, 0.42518887, 0.35281437, 0.27813674)
for ($i = 0; $i -lt $means.Length; $i++) {
    $ws.Cells.Item($i + 2, 52).Value = $means[$i]
}
